# Auto-generated cell updates applying the scheduled-runner market-data refresh.
# Each sheet corresponds to a crafting class; columns H-N hold cached market/profit data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 679.3570999999999
$ws.Range("I2").Value = 629.53845
$ws.Range("K2").Value = 629.53845
$ws.Range("M2").Value = -516.53845
$ws.Range("H11").Value = 174.5
$ws.Range("I11").Value = 174.5
$ws.Range("K11").Value = 174.5
$ws.Range("M11").Value = -34.5
$ws.Range("H43").Value = 5502.909
$ws.Range("I43").Value = 3880
$ws.Range("J43").Value = 6430.2856
$ws.Range("K43").Value = 3880
$ws.Range("L43").Value = 6430.2856
$ws.Range("M43").Value = -3811
$ws.Range("N43").Value = -6568.2856
$ws.Range("H46").Value = 127453.625
$ws.Range("J46").Value = 253750
$ws.Range("L46").Value = 761250
$ws.Range("N46").Value = -761488
$ws.Range("H60").Value = 127453.625
$ws.Range("J60").Value = 253750
$ws.Range("L60").Value = 761250
$ws.Range("N60").Value = -762218
$ws.Range("H76").Value = 4442.8335
$ws.Range("I76").Value = 4825
$ws.Range("J76").Value = 4251.75
$ws.Range("K76").Value = 4825
$ws.Range("L76").Value = 4251.75
$ws.Range("M76").Value = -4510
$ws.Range("N76").Value = -4881.75
$ws.Range("H79").Value = 4442.8335
$ws.Range("I79").Value = 4825
$ws.Range("J79").Value = 4251.75
$ws.Range("K79").Value = 4825
$ws.Range("L79").Value = 4251.75
$ws.Range("M79").Value = -3733
$ws.Range("N79").Value = -6435.75
$ws.Range("H132").Value = 5017.7666
$ws.Range("I132").Value = 4599.4614
$ws.Range("J132").Value = 7736.75
$ws.Range("K132").Value = 13798.3842
$ws.Range("L132").Value = 23210.25
$ws.Range("M132").Value = -11268.3842
$ws.Range("N132").Value = -28270.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4692.58
$ws.Range("J32").Value = 1196.25
$ws.Range("L32").Value = 1196.25
$ws.Range("N32").Value = -1770.25
$ws.Range("H46").Value = 8545.666999999999
$ws.Range("J46").Value = 8545.666999999999
$ws.Range("L46").Value = 8545.666999999999
$ws.Range("N46").Value = -9183.666999999999
$ws.Range("H63").Value = 3560
$ws.Range("I63").Value = 3755.7144
$ws.Range("J63").Value = 2875
$ws.Range("K63").Value = 3755.7144
$ws.Range("L63").Value = 2875
$ws.Range("M63").Value = -3069.7144
$ws.Range("N63").Value = -4247
$ws.Range("H66").Value = 3560
$ws.Range("I66").Value = 3755.7144
$ws.Range("J66").Value = 2875
$ws.Range("K66").Value = 18778.572
$ws.Range("L66").Value = 14375
$ws.Range("M66").Value = -15346.572
$ws.Range("N66").Value = -21239
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4070.55
$ws.Range("J105").Value = 6900
$ws.Range("L105").Value = 6900
$ws.Range("N105").Value = -10394
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 30586.6
$ws.Range("I132").Value = 9782.134
$ws.Range("K132").Value = 29346.402
$ws.Range("M132").Value = -26816.402
$ws.Range("H134").Value = 9721.875
$ws.Range("I134").Value = 9199.6
$ws.Range("K134").Value = 27598.8
$ws.Range("M134").Value = -25063.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 283.22223
$ws.Range("I10").Value = 193.625
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 580.875
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = -441.875
$ws.Range("N10").Value = -3278
$ws.Range("H57").Value = 5107.5
$ws.Range("I57").Value = 2929
$ws.Range("K57").Value = 8787
$ws.Range("M57").Value = -8228
$ws.Range("H75").Value = 533
$ws.Range("I75").Value = 299.5
$ws.Range("K75").Value = 898.5
$ws.Range("M75").Value = 99.5
$ws.Range("H78").Value = 533
$ws.Range("I78").Value = 299.5
$ws.Range("K78").Value = 2695.5
$ws.Range("M78").Value = 2296.5
$ws.Range("H110").Value = 23266.666
$ws.Range("I110").Value = 23266.666
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 69799.99800000001
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -65709.99800000001
$ws.Range("N110").ClearContents()
$ws.Range("H114").Value = 4332.6665
$ws.Range("I114").Value = 2499
$ws.Range("J114").Value = 4699.4
$ws.Range("K114").Value = 7497
$ws.Range("L114").Value = 14098.2
$ws.Range("M114").Value = -4243
$ws.Range("N114").Value = -20606.2
$ws.Range("H131").Value = 1949.0625
$ws.Range("I131").Value = 995
$ws.Range("J131").Value = 2085.3572
$ws.Range("K131").Value = 2985
$ws.Range("L131").Value = 6256.071599999999
$ws.Range("M131").Value = 2055
$ws.Range("N131").Value = -16336.0716
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 5000
$ws.Range("J59").Value = 2000
$ws.Range("L59").Value = 2000
$ws.Range("N59").Value = -3166
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H80").Value = 5018.5
$ws.Range("I80").Value = 3555.5
$ws.Range("J80").Value = 5750
$ws.Range("K80").Value = 3555.5
$ws.Range("L80").Value = 5750
$ws.Range("M80").Value = -2557.5
$ws.Range("N80").Value = -7746
$ws.Range("H83").Value = 5018.5
$ws.Range("I83").Value = 3555.5
$ws.Range("J83").Value = 5750
$ws.Range("K83").Value = 17777.5
$ws.Range("L83").Value = 28750
$ws.Range("M83").Value = -12785.5
$ws.Range("N83").Value = -38734
$ws.Range("H99").Value = 21399.6
$ws.Range("I99").Value = 19249.75
$ws.Range("K99").Value = 19249.75
$ws.Range("M99").Value = -17003.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6535.4287
$ws.Range("I100").Value = 7549.6
$ws.Range("K100").Value = 7549.6
$ws.Range("M100").Value = -7008.6
$ws.Range("H136").Value = 8977.032999999999
$ws.Range("I136").Value = 18363.428
$ws.Range("J136").Value = 6120.304
$ws.Range("K136").Value = 55090.284
$ws.Range("L136").Value = 18360.912
$ws.Range("M136").Value = -52540.284
$ws.Range("N136").Value = -23460.912
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 49999
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H96").Value = 3744.0625
$ws.Range("I96").Value = 3900.125
$ws.Range("J96").Value = 3588
$ws.Range("K96").Value = 3900.125
$ws.Range("L96").Value = 3588
$ws.Range("M96").Value = -2527.125
$ws.Range("N96").Value = -6334
$ws.Range("H132").Value = 19302.229
$ws.Range("I132").Value = 26232.137
$ws.Range("K132").Value = 78696.41099999999
$ws.Range("M132").Value = -76166.41099999999
$ws.Range("H136").Value = 1892.3889
$ws.Range("I136").Value = 1716.1818
$ws.Range("J136").Value = 2169.2856
$ws.Range("K136").Value = 5148.5454
$ws.Range("L136").Value = 6507.8568
$ws.Range("M136").Value = -2598.5454
$ws.Range("N136").Value = -11607.8568
